# Auto-generated COM-interop script to apply the WF 2017.xlsx diff
$wb = $excel.ActiveWorkbook
$chk = $wb.Worksheets.Item("Checking")
$sav = $wb.Worksheets.Item("Savings")

# ============================================================
# PHASE 1: Checking sheet - new rows 128-158 (numbers + formatting)
# ============================================================
$chk.Range("A3:D3").Copy() | Out-Null
$chk.Range("A128:D128").PasteSpecial(-4122) | Out-Null
$chk.Range("A128").Value = 42927
$chk.Range("B128").Value = -283
$chk.Range("C128").Clear() | Out-Null
$chk.Range("A3:D3").Copy() | Out-Null
$chk.Range("A129:D129").PasteSpecial(-4122) | Out-Null
$chk.Range("A129").Value = 42930
$chk.Range("B129").Value = -401.23
$chk.Range("C129").Clear() | Out-Null
$chk.Range("A3:D3").Copy() | Out-Null
$chk.Range("A130:D130").PasteSpecial(-4122) | Out-Null
$chk.Range("A130").Value = 42930
$chk.Range("B130").Value = 2236.81
$chk.Range("C130").Clear() | Out-Null
$chk.Range("A3:D3").Copy() | Out-Null
$chk.Range("A131:D131").PasteSpecial(-4122) | Out-Null
$chk.Range("A131").Value = 42933
$chk.Range("B131").Value = -1984.3
$chk.Range("C131").Clear() | Out-Null
$chk.Range("A3:D3").Copy() | Out-Null
$chk.Range("A132:D132").PasteSpecial(-4122) | Out-Null
$chk.Range("A132").Value = 42933
$chk.Range("B132").Value = 2242.41
$chk.Range("C132").Clear() | Out-Null
$chk.Range("A3:D3").Copy() | Out-Null
$chk.Range("A133:D133").PasteSpecial(-4122) | Out-Null
$chk.Range("A133").Value = 42934
$chk.Range("B133").Value = -1108
$chk.Range("C133").Value = 1767
$chk.Range("A3:D3").Copy() | Out-Null
$chk.Range("A134:D134").PasteSpecial(-4122) | Out-Null
$chk.Range("A134").Value = 42934
$chk.Range("B134").Value = -174.57
$chk.Range("C134").Clear() | Out-Null
$chk.Range("A3:D3").Copy() | Out-Null
$chk.Range("A135:D135").PasteSpecial(-4122) | Out-Null
$chk.Range("A135").Value = 42937
$chk.Range("B135").Value = -87.24
$chk.Range("C135").Clear() | Out-Null
$chk.Range("A3:D3").Copy() | Out-Null
$chk.Range("A136:D136").PasteSpecial(-4122) | Out-Null
$chk.Range("A136").Value = 42940
$chk.Range("B136").Value = -35.44
$chk.Range("C136").Value = 1768
$chk.Range("A3:D3").Copy() | Out-Null
$chk.Range("A137:D137").PasteSpecial(-4122) | Out-Null
$chk.Range("A137").Value = 42940
$chk.Range("B137").Value = -24.92
$chk.Range("C137").Value = 1771
$chk.Range("A2:D2").Copy() | Out-Null
$chk.Range("A138:D138").PasteSpecial(-4122) | Out-Null
$chk.Range("A138").Value = 42940
$chk.Range("B138").Value = -29.11
$chk.Range("C138").Clear() | Out-Null
$chk.Range("A3:D3").Copy() | Out-Null
$chk.Range("A139:D139").PasteSpecial(-4122) | Out-Null
$chk.Range("A139").Value = 42941
$chk.Range("B139").Value = 25.96
$chk.Range("C139").Clear() | Out-Null
$chk.Range("A3:D3").Copy() | Out-Null
$chk.Range("A140:D140").PasteSpecial(-4122) | Out-Null
$chk.Range("A140").Value = 42944
$chk.Range("B140").Value = 2236.8200000000002
$chk.Range("C140").Clear() | Out-Null
$chk.Range("A3:D3").Copy() | Out-Null
$chk.Range("A141:D141").PasteSpecial(-4122) | Out-Null
$chk.Range("A141").Value = 42947
$chk.Range("B141").Value = 2.17
$chk.Range("C141").Clear() | Out-Null
$chk.Range("A3:D3").Copy() | Out-Null
$chk.Range("A142:D142").PasteSpecial(-4122) | Out-Null
$chk.Range("A142").Value = 42947
$chk.Range("B142").Value = -33.130000000000003
$chk.Range("C142").Clear() | Out-Null
$chk.Range("A3:D3").Copy() | Out-Null
$chk.Range("A143:D143").PasteSpecial(-4122) | Out-Null
$chk.Range("A143").Value = 42947
$chk.Range("B143").Value = -85.94
$chk.Range("C143").Clear() | Out-Null
$chk.Range("A2:D2").Copy() | Out-Null
$chk.Range("A144:D144").PasteSpecial(-4122) | Out-Null
$chk.Range("A144").Value = 42947
$chk.Range("B144").Value = -57.29
$chk.Range("C144").Clear() | Out-Null
$chk.Range("A3:D3").Copy() | Out-Null
$chk.Range("A145:D145").PasteSpecial(-4122) | Out-Null
$chk.Range("A145").Value = 42948
$chk.Range("B145").Value = -75
$chk.Range("C145").Clear() | Out-Null
$chk.Range("A3:D3").Copy() | Out-Null
$chk.Range("A146:D146").PasteSpecial(-4122) | Out-Null
$chk.Range("A146").Value = 42954
$chk.Range("B146").Value = -1907.05
$chk.Range("C146").Clear() | Out-Null
$chk.Range("A3:D3").Copy() | Out-Null
$chk.Range("A147:D147").PasteSpecial(-4122) | Out-Null
$chk.Range("A147").Value = 42954
$chk.Range("B147").Value = -393
$chk.Range("C147").Clear() | Out-Null
$chk.Range("A2:D2").Copy() | Out-Null
$chk.Range("A148:D148").PasteSpecial(-4122) | Out-Null
$chk.Range("A148").Value = 42954
$chk.Range("B148").Value = -596.67999999999995
$chk.Range("C148").Clear() | Out-Null
$chk.Range("A3:D3").Copy() | Out-Null
$chk.Range("A149:D149").PasteSpecial(-4122) | Out-Null
$chk.Range("A149").Value = 42954
$chk.Range("B149").Value = -253.8
$chk.Range("C149").Value = 1772
$chk.Range("A2:D2").Copy() | Out-Null
$chk.Range("A150:D150").PasteSpecial(-4122) | Out-Null
$chk.Range("A150").Value = 42955
$chk.Range("B150").Value = -95
$chk.Range("C150").Clear() | Out-Null
$chk.Range("A3:D3").Copy() | Out-Null
$chk.Range("A151:D151").PasteSpecial(-4122) | Out-Null
$chk.Range("A151").Value = 42956
$chk.Range("B151").Value = 3195
$chk.Range("C151").Clear() | Out-Null
$chk.Range("A3:D3").Copy() | Out-Null
$chk.Range("A152:D152").PasteSpecial(-4122) | Out-Null
$chk.Range("A152").Value = 42958
$chk.Range("B152").Value = 2236.8200000000002
$chk.Range("C152").Clear() | Out-Null
$chk.Range("A3:D3").Copy() | Out-Null
$chk.Range("A153:D153").PasteSpecial(-4122) | Out-Null
$chk.Range("A153").Value = 42961
$chk.Range("B153").Value = 241.97
$chk.Range("C153").Clear() | Out-Null
$chk.Range("A2:D2").Copy() | Out-Null
$chk.Range("A154:D154").PasteSpecial(-4122) | Out-Null
$chk.Range("A154").Value = 42963
$chk.Range("B154").Value = 9900
$chk.Range("C154").Clear() | Out-Null
$chk.Range("A3:D3").Copy() | Out-Null
$chk.Range("A155:D155").PasteSpecial(-4122) | Out-Null
$chk.Range("A155").Value = 42964
$chk.Range("B155").Value = -1575.9
$chk.Range("C155").Clear() | Out-Null
$chk.Range("A3:D3").Copy() | Out-Null
$chk.Range("A156:D156").PasteSpecial(-4122) | Out-Null
$chk.Range("A156").Value = 42965
$chk.Range("B156").Value = -135.58000000000001
$chk.Range("C156").Value = 1773
$chk.Range("A2:D2").Copy() | Out-Null
$chk.Range("A157:D157").PasteSpecial(-4122) | Out-Null
$chk.Range("A157").Value = 42968
$chk.Range("B157").Value = 3884.02
$chk.Range("C157").Clear() | Out-Null
$chk.Range("A3:D3").Copy() | Out-Null
$chk.Range("A158:D158").PasteSpecial(-4122) | Out-Null
$chk.Range("A158").Value = 42972
$chk.Range("B158").Value = 2236.8200000000002
$chk.Range("C158").Clear() | Out-Null

# ============================================================
# PHASE 2: Savings sheet - new rows 21 and 23 (numbers + formatting)
# (row 22 stays fully blank, matching the source statement layout)
# ============================================================
$sav.Range("A12:C12").Copy() | Out-Null
$sav.Range("A21:C21").PasteSpecial(-4122) | Out-Null
$sav.Range("A21").Value = 42947
$sav.Range("B21").Value = 11.05
$sav.Range("A12:C12").Copy() | Out-Null
$sav.Range("A23:C23").PasteSpecial(-4122) | Out-Null
$sav.Range("A23").Value = 42948
$sav.Range("B23").Value = 75
# column E mirrors column C formatting (s=5 / fontId 3, no number format)
$sav.Range("C21").Copy() | Out-Null
$sav.Range("E21").PasteSpecial(-4122) | Out-Null
$sav.Range("C23").Copy() | Out-Null
$sav.Range("E23").PasteSpecial(-4122) | Out-Null

# ============================================================
# PHASE 3: shared-string text, written in the exact sequence the
# source workbook first introduced each string (keeps the
# xl/sharedStrings.xml table byte-identical to the target).
# ============================================================
$sav.Range("C21").Value = "*"
$sav.Range("E23").Value = "RECURRING TRANSFER FROM KLEINFINGHER L WELLS FARGO PRIME CHECKING REF #OP03MMX6KR XXXXXX1140"
$chk.Range("D128").Value = "ExpressSend Transaction wf386952602 Banco Santander MXN 5000 FX Rate: 17.6677 USD Amt: 283.00 To: 60577934998 LILIAN PAULINA KLEINFINGHER LECUONA"
$chk.Range("D130").Value = "AMGEN INC DIRECT DEP 170714 942303672515VH1 KLEINFINGHER L,LIZ Y"
$chk.Range("D129").Value = "ExpressSend Transaction wf160808098 Banco Santander MXN 7000 FX Rate: 17.4463 USD Amt: 401.23 To: 60577934998 LILIAN PAULINA KLEINFINGHER LECUONA"
$chk.Range("D132").Value = "ATM CHECK DEPOSIT ON 07/17 1596 NORTH MOORPARK R THOUSAND OAKS CA 0000268 ATM ID 9842E CARD 9001"
$chk.Range("D131").Value = "AMERICAN EXPRESS ACH PMT 170717 A8606 LIZ KLEINFINGHER"
$chk.Range("D134").Value = "ExpressSend Transaction wf189486160 Banco Santander MXN 3000 FX Rate: 17.185 USD Amt: 174.57 To: 60577934998 LILIAN PAULINA KLEINFINGHER LECUONA"
$chk.Range("D133").Value = "CHECK # 1767"
$chk.Range("D135").Value = "ExpressSend Transaction wf873807465 Banco Santander MXN 1500 FX Rate: 17.1948 USD Amt: 87.24 To: 60577934998 LILIAN PAULINA KLEINFINGHER LECUONA"
$chk.Range("D138").Value = "ExpressSend Transaction wf957085873 Banco Santander MXN 500 FX Rate: 17.1751 USD Amt: 29.11 To: 60577934998 LILIAN PAULINA KLEINFINGHER LECUONA"
$chk.Range("D137").Value = "CHECK # 1771"
$chk.Range("D136").Value = "CHECK # 1768"
$chk.Range("D139").Value = "MOBILE DEPOSIT : REF NUMBER :718250737277"
$chk.Range("D140").Value = "AMGEN INC DIRECT DEP 170728 932903770633VH1 KLEINFINGHER L,LIZ Y"
$chk.Range("D144").Value = "ExpressSend Transaction wf280766470 Banco Santander MXN 1000 FX Rate: 17.455 USD Amt: 57.29 To: 60577934998 LILIAN PAULINA KLEINFINGHER LECUONA"
$chk.Range("D143").Value = "ExpressSend Transaction wf238447900 Soriana- Comercial Mexicana MXN 1500 FX Rate: 17.455 USD Amt: 85.94 To: CASH PICKUP CARLOS ARTURO MORA LECUONA"
$chk.Range("D142").Value = "Check # 1770 (Converted ACH) SOCALGAS ARC PYMT 170728 1770 0106108548"
$chk.Range("D145").Value = "RECURRING TRANSFER TO KLEINFINGHER L SAVINGS REF #OP03MMX6KR XXXXXX2638"
$chk.Range("D149").Value = "CHECK # 1772"
$chk.Range("D148").Value = "CHASE CREDIT CRD AUTOPAY 170806 000000000053157 KLEINFINGHER LIZ"
$chk.Range("D147").Value = "OAK PARK VILLAGE PAYMENTS 170807 00143-2936 KLEINFINGHER, LIZ Y"
$chk.Range("D146").Value = "WF HOME MTG AUTO PAY 080517 0337643209 LIZ Y KLEINFINGHER"
$chk.Range("D150").Value = "CITI AUTOPAY PAYMENT 170807 082402802527302 LIZ KLEINFINGHER"
$chk.Range("D151").Value = "ATM CHECK DEPOSIT ON 08/09 220 N MOORPARK RD THOUSAND OAKS CA 0008026 ATM ID 9849A CARD 9001"
$chk.Range("D152").Value = "AMGEN INC DIRECT DEP 170811 675048651278VH1 KLEINFINGHER L,LIZ Y"
$chk.Range("D153").Value = "MOBILE DEPOSIT : REF NUMBER :115130539452"
$chk.Range("D154").Value = "WT 1132737053015208 BBVA BANCOMER SA /ORG=1/LILIAN PAULINA KLEINFINGHER LE SRF# 1132737053015208 TRN#170816111390 RFB#"
$chk.Range("D155").Value = "AMERICAN EXPRESS ACH PMT 170817 A1730 LIZ KLEINFINGHER"
$chk.Range("D156").Value = "CHECK # 1773"
$chk.Range("D157").Value = "WT 1132774519010431 BBVA BANCOMER SA /ORG=1/LILIAN PAULINA KLEINFINGHER LE SRF# 1132774519010431 TRN#170818130129 RFB#"
$chk.Range("D158").Value = "AMGEN INC DIRECT DEP 170825 611042582835VH1 KLEINFINGHER L,LIZ Y"

# Reused existing strings (no new shared-string slot needed)
$chk.Range("D141").Value = "INTEREST PAYMENT"
$sav.Range("E21").Value = "INTEREST PAYMENT"

# ============================================================
# PHASE 4: dimension / view / sort-state bookkeeping
# ============================================================
